$wb = $excel.ActiveWorkbook

# The source sheet being brought in here originally came from a workbook
# saved with the legacy Excel default style (Arial 10pt "Normal"); bringing
# it in registers that style alongside the destination workbook's own
# "Normal" style, hence "Normal 2".
$style = $wb.Styles.Add("Normal 2")
$style.Font.Size = 10
$style.Font.Name = "Arial"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "mi_column_empty_rows"

$ws.Cells.Item(1, 1).Value = "a"
$ws.Cells.Item(1, 2).Value = "b"
$ws.Cells.Item(2, 1).Value = "A"
$ws.Cells.Item(2, 2).Value = "B"
$ws.Cells.Item(5, 1).Value = 1
$ws.Cells.Item(5, 2).Value = 3
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 2).Value = 4
